# Update the date heading.
$d = $word.ActiveDocument
$p = $d.Paragraphs.Item(1)
$p.Range.Text = "2024-07-20 Saturday"

# Update the division-problem answers in the practice table.
# Only the data rows (1, 5, 9, 13, 17) of the 20-row / 5-column table
# carry text; the intervening rows are blank spacer rows.
$t = $d.Tables.Item(1)

$t.Cell(1,1).Range.Text  = "18÷7=2, 4"
$t.Cell(1,2).Range.Text  = "21÷8=2, 5"
$t.Cell(1,3).Range.Text  = "15÷3=5, 0"
$t.Cell(1,4).Range.Text  = "79÷6=13, 1"
$t.Cell(1,5).Range.Text  = "36÷9=4, 0"

$t.Cell(5,1).Range.Text  = "36÷9=4, 0"
$t.Cell(5,2).Range.Text  = "73÷5=14, 3"
$t.Cell(5,3).Range.Text  = "39÷9=4, 3"
$t.Cell(5,4).Range.Text  = "80÷5=16, 0"
$t.Cell(5,5).Range.Text  = "36÷7=5, 1"

$t.Cell(9,1).Range.Text  = "39÷6=6, 3"
$t.Cell(9,2).Range.Text  = "71÷3=23, 2"
$t.Cell(9,3).Range.Text  = "86÷4=21, 2"
$t.Cell(9,4).Range.Text  = "55÷3=18, 1"
$t.Cell(9,5).Range.Text  = "85÷8=10, 5"

$t.Cell(13,1).Range.Text = "40÷3=13, 1"
$t.Cell(13,2).Range.Text = "78÷9=8, 6"
$t.Cell(13,3).Range.Text = "87÷3=29, 0"
$t.Cell(13,4).Range.Text = "72÷4=18, 0"
$t.Cell(13,5).Range.Text = "60÷5=12, 0"

$t.Cell(17,1).Range.Text = "61÷5=12, 1"
$t.Cell(17,2).Range.Text = "68÷4=17, 0"
$t.Cell(17,3).Range.Text = "29÷2=14, 1"
$t.Cell(17,4).Range.Text = "81÷9=9, 0"
$t.Cell(17,5).Range.Text = "67÷3=22, 1"
